# Refresh cryptos list: update Price (D) and Volume(1h) (E) columns per latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.378.09"
$ws.Range("E2").Value = "  +2.56%  "

$ws.Range("D3").Value = "1.824.77"
$ws.Range("E3").Value = "  +1.71%  "

$ws.Range("D4").Value = "'1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.13%  "

$ws.Range("D5").Value = "'313.93"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.55%  "

$ws.Range("D6").Value = "'1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.09%  "

$ws.Range("D7").Value = "'0.4672"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.81%  "

$ws.Range("E8").Value = "  +3.52%  "

$ws.Range("D9").Value = "'0.07438"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.98%  "

$ws.Range("D10").Value = "'0.8752"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.38%  "

$ws.Range("D11").Value = "'20.78"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.10%  "

$ws.Range("D12").Value = "1.824.77"
$ws.Range("E12").Value = "  -3.20%  "

$ws.Range("D13").Value = "'6.688"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.35%  "

$ws.Range("D14").Value = "'5.418"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.89%  "

$ws.Range("D15").Value = "'92.98"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.13%  "

$ws.Range("D16").Value = "'0.07087"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.21%  "

$ws.Range("E17").Value = "  -0.09%  "

$ws.Range("D18").Value = "'0.000008795"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.48%  "

$ws.Range("D19").Value = "'0.9999"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.07%  "

$ws.Range("D20").Value = "'15.02"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.49%  "

$ws.Range("D21").Value = "27.371.69"
$ws.Range("E21").Value = "  +2.28%  "

$ws.Range("E22").Value = "  +3.50%  "

$ws.Range("E23").Value = "  +2.19%  "

$ws.Range("D24").Value = "2.052.52"
$ws.Range("E24").Value = "  -4.52%  "

$ws.Range("D25").Value = "'1.939"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.28%  "

$ws.Range("D26").Value = "'151.18"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.47%  "

$ws.Range("D27").Value = "'2.251"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.84%  "

$ws.Range("D28").Value = "'18.66"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.31%  "

$ws.Range("D29").Value = "'5.337"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.15%  "

$ws.Range("D30").Value = "'117.13"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.74%  "

$ws.Range("D31").Value = "'0.08954"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.91%  "

$ws.Range("D32").Value = "'0.7876"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.43%  "

$ws.Range("D33").Value = "'1.193"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.53%  "

$ws.Range("D34").Value = "'4.537"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.57%  "

$ws.Range("E35").Value = "  +0.14%  "

$ws.Range("D36").Value = "'1.000"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.09%  "

$ws.Range("D37").Value = "'1.102"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.70%  "

$ws.Range("D38").Value = "'0.01975"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.03%  "

$ws.Range("D39").Value = "'0.05246"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.50%  "

$ws.Range("E40").Value = "  +3.91%  "

$ws.Range("D41").Value = "'0.5367"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.88%  "

$ws.Range("D42").Value = "'2.900"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.95%  "

$ws.Range("D43").Value = "'2.351"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +20.16%  "

$ws.Range("D44").Value = "'0.1704"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.66%  "

$ws.Range("D45").Value = "'8.656"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.40%  "

$ws.Range("D46").Value = "'0.5097"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.25%  "

$ws.Range("D47").Value = "'10.64"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.92%  "

$ws.Range("D48").Value = "'105.92"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.47%  "

$ws.Range("E49").Value = "  +1.67%  "

$ws.Range("D50").Value = "'0.9997"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.05%  "

$ws.Range("E51").Value = "  +1.41%  "

